$wb = $excel.ActiveWorkbook

# Sheet "high_loadings" - update Category (column B) values
$ws1 = $wb.Worksheets.Item("high_loadings")
$ws1.Range("B4").Value = 2
$ws1.Range("B6").Value = 2
$ws1.Range("B7").Value = 3
$ws1.Range("B9").Value = 2
$ws1.Range("B11").Value = 1
$ws1.Range("B13").Value = 1
$ws1.Range("B14").Value = 1
$ws1.Range("B15").Value = 1
$ws1.Range("B16").Value = 2
$ws1.Range("B17").Value = 2
$ws1.Range("B18").Value = 2
$ws1.Range("B19").Value = 3
$ws1.Range("B22").Value = 3
$ws1.Range("B23").Value = 2
$ws1.Range("B24").Value = 1
$ws1.Range("B25").Value = 2
$ws1.Range("B27").Value = 3
$ws1.Range("B28").Value = 1
$ws1.Range("B29").Value = 2
$ws1.Range("B30").Value = 2
$ws1.Range("B31").Value = 3
$ws1.Range("B33").Value = 2
$ws1.Range("B34").Value = 3
$ws1.Range("B36").Value = 2
$ws1.Range("B37").Value = 3
$ws1.Range("B39").Value = 1
$ws1.Range("B40").Value = 3
$ws1.Range("B41").Value = 1
$ws1.Range("B42").Value = 3
$ws1.Range("B45").Value = 2
$ws1.Range("B47").Value = 3
$ws1.Range("B50").Value = 1
$ws1.Range("B52").Value = 1
$ws1.Range("B53").Value = 1
$ws1.Range("B55").Value = 1
$ws1.Range("B58").Value = 2
$ws1.Range("B60").Value = 2
$ws1.Range("B61").Value = 2
$ws1.Range("B62").Value = 1
$ws1.Range("B63").Value = 3
$ws1.Range("B65").Value = 1
$ws1.Range("B67").Value = 3
$ws1.Range("B69").Value = 2
$ws1.Range("B70").Value = 2
$ws1.Range("B71").Value = 3

# Sheet "Influencers_uniques" - update Category (column B) values
$ws2 = $wb.Worksheets.Item("Influencers_uniques")
$ws2.Range("B2").Value = 2
$ws2.Range("B3").Value = 3
$ws2.Range("B4").Value = 2
$ws2.Range("B6").Value = 1
$ws2.Range("B7").Value = 1
$ws2.Range("B8").Value = 2
$ws2.Range("B9").Value = 2
$ws2.Range("B10").Value = 3
$ws2.Range("B11").Value = 3
$ws2.Range("B15").Value = 2
$ws2.Range("B17").Value = 1
$ws2.Range("B18").Value = 2
$ws2.Range("B19").Value = 3
$ws2.Range("B20").Value = 3
$ws2.Range("B21").Value = 3
$ws2.Range("B23").Value = 2
$ws2.Range("B24").Value = 2
$ws2.Range("B25").Value = 1
$ws2.Range("B27").Value = 2
$ws2.Range("B30").Value = 3
$ws2.Range("B31").Value = 1
$ws2.Range("B32").Value = 3
$ws2.Range("B34").Value = 3
$ws2.Range("B35").Value = 2
$ws2.Range("B36").Value = 2
$ws2.Range("B39").Value = 1
$ws2.Range("B41").Value = 3
$ws2.Range("B43").Value = 2
$ws2.Range("B45").Value = 1
$ws2.Range("B46").Value = 2
$ws2.Range("B47").Value = 1
$ws2.Range("B49").Value = 1
$ws2.Range("B50").Value = 1
$ws2.Range("B51").Value = 3
$ws2.Range("B52").Value = 2
$ws2.Range("B53").Value = 2
$ws2.Range("B54").Value = 3
$ws2.Range("B55").Value = 3
$ws2.Range("B57").Value = 1
$ws2.Range("B61").Value = 3
$ws2.Range("B62").Value = 3
$ws2.Range("B63").Value = 1
$ws2.Range("B64").Value = 3
$ws2.Range("B65").Value = 2
$ws2.Range("B66").Value = 2
$ws2.Range("B67").Value = 2
$ws2.Range("B68").Value = 1
$ws2.Range("B69").Value = 3
$ws2.Range("B70").Value = 2
$ws2.Range("B71").Value = 2
